$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1237.9667
$ws.Range("J17").Value = 1372.4584
$ws.Range("L17").Value = 4117.3752
$ws.Range("N17").Value = -4453.3752
$ws.Range("H40").Value = 4000
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H74").Value = 6979.1665
$ws.Range("I74").Value = 6979.1665
$ws.Range("K74").Value = 6979.1665
$ws.Range("M74").Value = -6043.1665
$ws.Range("H77").Value = 6979.1665
$ws.Range("I77").Value = 6979.1665
$ws.Range("K77").Value = 34895.8325
$ws.Range("M77").Value = -30215.8325
$ws.Range("H112").Value = 2451.6155
$ws.Range("J112").Value = 2999
$ws.Range("L112").Value = 8997
$ws.Range("N112").Value = -11213
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1880.5555
$ws.Range("I2").Value = 1860.7142
$ws.Range("J2").Value = 1950
$ws.Range("K2").Value = 1860.7142
$ws.Range("L2").Value = 1950
$ws.Range("M2").Value = -1747.7142
$ws.Range("N2").Value = -2176
$ws.Range("H32").Value = 11723.444
$ws.Range("I32").Value = 11060.117
$ws.Range("K32").Value = 11060.117
$ws.Range("M32").Value = -10773.117
$ws.Range("H45").Value = 36664.168
$ws.Range("I45").Value = 9996.25
$ws.Range("J45").Value = 90000
$ws.Range("K45").Value = 9996.25
$ws.Range("L45").Value = 90000
$ws.Range("M45").Value = -9619.25
$ws.Range("N45").Value = -90754
$ws.Range("H55").Value = 60053
$ws.Range("J55").Value = 60053
$ws.Range("L55").Value = 60053
$ws.Range("N55").Value = -60683
$ws.Range("H61").Value = 4155.875
$ws.Range("I61").Value = 4035.2856
$ws.Range("K61").Value = 4035.2856
$ws.Range("M61").Value = -3823.2856
$ws.Range("H63").Value = 6060.143
$ws.Range("J63").Value = 10496.667
$ws.Range("L63").Value = 10496.667
$ws.Range("N63").Value = -11868.667
$ws.Range("H66").Value = 6060.143
$ws.Range("J66").Value = 10496.667
$ws.Range("L66").Value = 52483.335
$ws.Range("N66").Value = -59347.335
$ws.Range("H74").Value = 8128
$ws.Range("I74").Value = 10004
$ws.Range("K74").Value = 10004
$ws.Range("M74").Value = -9130
$ws.Range("H77").Value = 8128
$ws.Range("I77").Value = 10004
$ws.Range("K77").Value = 50020
$ws.Range("M77").Value = -45652
$ws.Range("H102").Value = 2065.4
$ws.Range("I102").Value = 2065.4
$ws.Range("K102").Value = 2065.4
$ws.Range("M102").Value = -443.4000000000001
$ws.Range("H116").Value = 1880.5555
$ws.Range("I116").Value = 1860.7142
$ws.Range("J116").Value = 1950
$ws.Range("K116").Value = 1860.7142
$ws.Range("L116").Value = 1950
$ws.Range("M116").Value = 433.2858000000001
$ws.Range("N116").Value = -6538
$ws.Range("H136").Value = 4155.875
$ws.Range("I136").Value = 4035.2856
$ws.Range("K136").Value = 12105.8568
$ws.Range("M136").Value = -9555.856800000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1880.5555
$ws.Range("I3").Value = 1860.7142
$ws.Range("J3").Value = 1950
$ws.Range("K3").Value = 1860.7142
$ws.Range("L3").Value = 1950
$ws.Range("M3").Value = -1746.7142
$ws.Range("N3").Value = -2178

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3153.5833
$ws.Range("I58").Value = 3242.889
$ws.Range("K58").Value = 3242.889
$ws.Range("M58").Value = -3039.889
$ws.Range("H122").Value = 1198.6364
$ws.Range("J122").Value = 803
$ws.Range("L122").Value = 2409
$ws.Range("N122").Value = -7309
$ws.Range("H132").Value = 147644.58
$ws.Range("I132").Value = 502000
$ws.Range("J132").Value = 5902.4
$ws.Range("K132").Value = 1506000
$ws.Range("L132").Value = 17707.2
$ws.Range("M132").Value = -1503470
$ws.Range("N132").Value = -22767.2
$ws.Range("H134").Value = 9564.4
$ws.Range("I134").Value = 9564.4
$ws.Range("K134").Value = 28693.2
$ws.Range("M134").Value = -26158.2
$ws.Range("H136").Value = 3153.5833
$ws.Range("I136").Value = 3242.889
$ws.Range("K136").Value = 9728.667000000001
$ws.Range("M136").Value = -7178.667000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 450
$ws.Range("I86").Value = 450
$ws.Range("K86").Value = 1350
$ws.Range("M86").Value = -164
$ws.Range("H89").Value = 450
$ws.Range("I89").Value = 450
$ws.Range("K89").Value = 4050
$ws.Range("M89").Value = 1878
$ws.Range("H97").Value = 522.05
$ws.Range("I97").Value = 496.89474
$ws.Range("K97").Value = 1490.68422
$ws.Range("M97").Value = -994.6842200000001
$ws.Range("H113").Value = 2119.6667
$ws.Range("J113").Value = 2343.6
$ws.Range("L113").Value = 7030.799999999999
$ws.Range("N113").Value = -11370.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 7333
$ws.Range("I27").Value = 6999
$ws.Range("J27").Value = 7500
$ws.Range("K27").Value = 6999
$ws.Range("L27").Value = 7500
$ws.Range("M27").Value = -6833
$ws.Range("N27").Value = -7832
$ws.Range("H132").Value = 4411.909
$ws.Range("I132").Value = 3576.1428
$ws.Range("K132").Value = 10728.4284
$ws.Range("M132").Value = -8198.428400000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3500
$ws.Range("I46").Value = 3500
$ws.Range("K46").Value = 3500
$ws.Range("M46").Value = -3312
$ws.Range("H55").Value = 566.8570999999999
$ws.Range("I55").Value = 648.9091
$ws.Range("J55").Value = 266
$ws.Range("K55").Value = 648.9091
$ws.Range("L55").Value = 266
$ws.Range("M55").Value = -475.9091
$ws.Range("N55").Value = -612
$ws.Range("H136").Value = 5845.0835
$ws.Range("I136").Value = 2960.111
$ws.Range("K136").Value = 8880.332999999999
$ws.Range("M136").Value = -6330.332999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 21880
$ws.Range("J54").Value = 21880
$ws.Range("L54").Value = 21880
$ws.Range("N54").Value = -22920
$ws.Range("H70").Value = 70000
$ws.Range("J70").Value = 70000
$ws.Range("L70").Value = 70000
$ws.Range("N70").Value = -70630
$ws.Range("H73").Value = 70000
$ws.Range("J73").Value = 70000
$ws.Range("L73").Value = 70000
$ws.Range("N73").Value = -72184
$ws.Range("H132").Value = 2747.8333
$ws.Range("I132").Value = 1622.25
$ws.Range("K132").Value = 4866.75
$ws.Range("M132").Value = -2336.75
